$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "30.333.61"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(2, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.18%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(3, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.009.11"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(3, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.87%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(4, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.16%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "324.76"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(5, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.45%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(7, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5130"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(7, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.72%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(8, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4259"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(8, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.37%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(9, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08705"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(9, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.05%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(10, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.87%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(11, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "43.17"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(11, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.82%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.72"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(12, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.08%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(13, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.009.24"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(13, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.04%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(14, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.571"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(14, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.73%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(15, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.472"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(15, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.43%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(16, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(17, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "94.36"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(17, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.29%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(18, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.72%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(19, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06542"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(19, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.59%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(20, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.85"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(20, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.78%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(21, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(22, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.44%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "30.384.36"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(23, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.25%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(24, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.51%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.258"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(25, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.90%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.245.78"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(26, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.24%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(27, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.41"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(27, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.16%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(28, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.12%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.418"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(29, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.17%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(30, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.59%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(31, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.56%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1053"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(32, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.69%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.078"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.30%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(34, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.827"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(34, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(35, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.369"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(35, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +14.52%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(36, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.11%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(37, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06682"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(37, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.32%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(38, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.460"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(38, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.32%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(39, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.36"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(39, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +8.54%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(40, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.102"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(40, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.23%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(41, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2192"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(41, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.78%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(42, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6636"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(42, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.64%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(44, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.14%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.60"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.96%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(46, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6167"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(46, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.98%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.184"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(47, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.44%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(48, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.662"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(48, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.67%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(49, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.62%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "124.19"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(50, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.63%  "
$cell.Style = $origStyle

$cell = $ws.Cells.Item(51, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "80.69"
$cell.Style = $origStyle

$cell = $ws.Cells.Item(51, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.28%  "
$cell.Style = $origStyle

